$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Rename the sheet "Sheet 1" -> "Etapa"
# ---------------------------------------------------------------------------
$ws.Name = "Etapa"

# ---------------------------------------------------------------------------
# 2. View tweaks: hide gridlines, freeze the header row
# ---------------------------------------------------------------------------
$ws.Activate()
$excel.ActiveWindow.DisplayGridlines = $false
$ws.Range("A2").Select() | Out-Null
$excel.ActiveWindow.FreezePanes = $true

# ---------------------------------------------------------------------------
# 3. Column widths (approximate character widths from the template)
# ---------------------------------------------------------------------------
$colWidths = @{
    1  = 12.75   # A  etapa
    2  = 6.75    # B  n_total
    3  = 12.75   # C  n_error_fecha
    4  = 14.75   # D  pct_error_fecha
    5  = 16.75   # E  n_regla_operativa
    6  = 18.75   # F  pct_regla_operativa
    7  = 8.75    # G  n_validos
    8  = 9.75    # H  n_outliers
    9  = 11.75   # I  pct_outliers
    10 = 21.75   # J  n_validos_sin_outliers
    11 = 7.75    # K  media
    12 = 7.75    # L  mediana
    13 = 7.75    # M  p25
    14 = 7.75    # N  p75
    15 = 7.75    # O  minimo
    16 = 7.75    # P  maximo
}
foreach ($col in $colWidths.Keys) {
    $ws.Columns.Item($col).ColumnWidth = $colWidths[$col]
}

# ---------------------------------------------------------------------------
# 4. Header row styling (row 1): white bold Calibri on dark blue fill,
#    thin bottom border, centered both ways
# ---------------------------------------------------------------------------
$hdr = $ws.Range("A1:P1")
$hdr.Font.Name = "Calibri"
$hdr.Font.Color = 0xFFFFFF
$hdr.Font.Bold = $true
$hdr.Interior.Color = 0x794E1F
$hdr.Borders.Item(9).LineStyle = 1
$hdr.Borders.Item(9).Color = 0x000000
$hdr.HorizontalAlignment = -4108
$hdr.VerticalAlignment = -4108

# ---------------------------------------------------------------------------
# 5. Body rows (2-4): black Calibri text, boxed border on every cell,
#    vertically centered
# ---------------------------------------------------------------------------
for ($r = 2; $r -le 4; $r++) {
    for ($c = 1; $c -le 16; $c++) {
        $cell = $ws.Cells.Item($r, $c)
        $cell.Font.Name = "Calibri"
        $cell.Font.Color = 0x000000
        $cell.BorderAround(1, 2, 1, 0x000000)
        $cell.VerticalAlignment = -4108
    }
}

# Column A ("etapa"): left aligned, general format
$ws.Range("A2:A4").HorizontalAlignment = -4131

# Count columns (whole numbers): #,##0 format, centered
$countCols = @(2, 3, 5, 7, 8, 10)
foreach ($c in $countCols) {
    $rng = $ws.Range($ws.Cells.Item(2, $c), $ws.Cells.Item(4, $c))
    $rng.NumberFormat = "#,##0"
    $rng.HorizontalAlignment = -4108
}

# Percentage columns: 0.0"%" format, centered
$pctCols = @(4, 6, 9)
foreach ($c in $pctCols) {
    $rng = $ws.Range($ws.Cells.Item(2, $c), $ws.Cells.Item(4, $c))
    $rng.NumberFormat = "0.0""%"""
    $rng.HorizontalAlignment = -4108
}

# Duration columns (media..maximo): general format, centered
$ws.Range("K2:P4").HorizontalAlignment = -4108

# ---------------------------------------------------------------------------
# 6. AutoFilter across the full table + hidden _FilterDatabase defined name
# ---------------------------------------------------------------------------
$ws.Range("A1:P4").AutoFilter() | Out-Null
$fdb = $ws.Names.Add("_xlnm._FilterDatabase", "='Etapa'!`$A`$1:`$P`$4")
$fdb.Visible = $false

Write-Host "edit complete"
